$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91
$ws.Range("A91").Value = 90
$ws.Range("B91").Value = 11
$ws.Range("C91").Value = "Merchants Coin"
$ws.Range("D91").Value = "Merchants need their coin, don't let anyone get in the way of that!"
$ws.Range("E91").Value = 1
$ws.Range("I91").Value = 0.05
$ws.Range("J91").Value = 0.02
$ws.Range("N91").Value = 0.05

# Row 92
$ws.Range("A92").Value = 91
$ws.Range("B92").Value = 11
$ws.Range("C92").Value = "Merchants Sale"
$ws.Range("D92").Value = "What do you have for sale? Deals 1000 Damage growing by 10 over time for an additional 1000 damage using 5% of your damage stat."
$ws.Range("E92").Value = 12
$ws.Range("F92").Value = 1000
$ws.Range("G92").Value = 10
$ws.Range("H92").Value = 0.05
$ws.Range("I92").Value = 0.08
$ws.Range("N92").Value = 0.1
$ws.Range("O92").Value = "attack"

# Row 93
$ws.Range("A93").Value = 92
$ws.Range("B93").Value = 11
$ws.Range("C93").Value = "Merchants Charm"
$ws.Range("D93").Value = "Charm the enemy into letting their guard down. (Reduces enemy skills and resistances)"
$ws.Range("E93").Value = 24
$ws.Range("S93").Value = 0.15
$ws.Range("T93").Value = 0.1

# Row 94
$ws.Range("A94").Value = 93
$ws.Range("B94").Value = 11
$ws.Range("C94").Value = "Sales Pitch"
$ws.Range("D94").Value = "Make a sale to your enemy. They will do less damage now - through their affixes, but they might still kill you."
$ws.Range("E94").Value = 36
$ws.Range("Q94").Value = 0.25

# Row 95
$ws.Range("A95").Value = 94
$ws.Range("B95").Value = 11
$ws.Range("C95").Value = "Caravans War"
$ws.Range("D95").Value = "Lash out at the enemy with a caravans strength. Deal 10,000 Damage ground by 100 for an additional 10,000 damage. Use 25% of your damage stat."
$ws.Range("E95").Value = 48
$ws.Range("F95").Value = 10000
$ws.Range("G95").Value = 100
$ws.Range("H95").Value = 0.25
$ws.Range("I95").Value = 0.2
$ws.Range("J95").Value = 0.1
$ws.Range("K95").Value = 0.05
$ws.Range("N95").Value = 0.25
$ws.Range("O95").Value = "attack"
$ws.Range("Q95").Value = 0.1
$ws.Range("R95").Value = 0.05
$ws.Range("S95").Value = 0.15
$ws.Range("T95").Value = 0.1

# Row 96
$ws.Range("A96").Value = 95
$ws.Range("B96").Value = 11
$ws.Range("C96").Value = "Merchants Defence"
$ws.Range("D96").Value = "Defend your self child, the enemy comes! Deals 5,000 damage growing by 50 for an additional 5,000 damage over time using 15% of your damage but only when using defend."
$ws.Range("E96").Value = 60
$ws.Range("F96").Value = 5000
$ws.Range("G96").Value = 50
$ws.Range("H96").Value = 0.15
$ws.Range("I96").Value = 0.3
$ws.Range("J96").Value = 0.25
$ws.Range("K96").Value = 0.1
$ws.Range("M96").Value = 0.5
$ws.Range("N96").Value = 0.1
$ws.Range("O96").Value = "defend"
$ws.Range("P96").Value = 0.1
$ws.Range("Q96").Value = 0.2
$ws.Range("R96").Value = 0.1
$ws.Range("S96").Value = 0.1
$ws.Range("T96").Value = 0.2

# Row 97
$ws.Range("A97").Value = 96
$ws.Range("B97").Value = 11
$ws.Range("C97").Value = "Coin Flip"
$ws.Range("D97").Value = "Flip a coin and see what happens. Deals 25,000 damage growing by 250 over time for an additional 25,000 damage using 30% of your damage stat,"
$ws.Range("E97").Value = 70
$ws.Range("F97").Value = 25000
$ws.Range("G97").Value = 250
$ws.Range("H97").Value = 0.3
$ws.Range("I97").Value = 0.2
$ws.Range("N97").Value = 0.1
$ws.Range("O97").Value = "attack"

# Row 98
$ws.Range("A98").Value = 97
$ws.Range("B98").Value = 11
$ws.Range("C98").Value = "Caravans Last Stand"
$ws.Range("D98").Value = "traveling with a caravan of merchants allows you to call on them when you need them the most. Deals 50,000 damage growing by 500 over time, using 15% of your damage stat. This only works when defending."
$ws.Range("E98").Value = 80
$ws.Range("F98").Value = 50000
$ws.Range("G98").Value = 500
$ws.Range("H98").Value = 0.15
$ws.Range("I98").Value = 0.3
$ws.Range("J98").Value = 0.45
$ws.Range("K98").Value = 0.1
$ws.Range("M98").Value = 0.3
$ws.Range("N98").Value = 0.45
$ws.Range("O98").Value = "defend"
$ws.Range("P98").Value = 0.25
$ws.Range("Q98").Value = 0.25
$ws.Range("R98").Value = 0.25
$ws.Range("S98").Value = 0.25
$ws.Range("T98").Value = 0.25

# Row 99
$ws.Range("A99").Value = 98
$ws.Range("B99").Value = 11
$ws.Range("C99").Value = "Magical Trade"
$ws.Range("D99").Value = "Your dealings with the mages and the magical folk on the road have trained you well. Use your cast and attack to deal 60,000 damage growing by 600 over time for an additional 60,000 damage using 10% of your damage stat., The enemies spells are useless against you."
$ws.Range("E99").Value = 90
$ws.Range("F99").Value = 60000
$ws.Range("G99").Value = 600
$ws.Range("H99").Value = 0.1
$ws.Range("I99").Value = 0.1
$ws.Range("L99").Value = 0.3
$ws.Range("N99").Value = 0.05
$ws.Range("O99").Value = "attack_and_cast"
$ws.Range("P99").Value = 0.75
$ws.Range("Q99").Value = 0.5
$ws.Range("R99").Value = 0.25
$ws.Range("S99").Value = 0.1
$ws.Range("T99").Value = 0.2
